$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing header row (A1:G1) holds: suggestedManufacturer, suggestedManufacturerPartNumber,
# unitPrice, totalPrice, leadTime, shippingTerms, comment.
# We need to shift these 7 existing columns to H1:N1, and place 7 new headers in A1:G1:
# code, purchaseRequestNumber, shortText, quantity, uom, manufacturer, manufacturerPartNumber.
# Column width metadata (<cols>) must stay as-is (still referencing column indices 1,2,3,6,7),
# so we move the cell values directly rather than performing a true column insert/shift.

# Capture the current values of the first 7 header cells before overwriting anything.
$existingHeaders = @()
for ($col = 1; $col -le 7; $col++) {
    $existingHeaders += $ws.Cells.Item(1, $col).Value2
}

# Move the existing headers from columns A-G (1-7) to columns H-N (8-14).
for ($col = 1; $col -le 7; $col++) {
    $ws.Cells.Item(1, $col + 7).Value2 = $existingHeaders[$col - 1]
}

# Write the new headers into columns A-G (1-7).
$newHeaders = @("code", "purchaseRequestNumber", "shortText", "quantity", "uom", "manufacturer", "manufacturerPartNumber")
for ($col = 1; $col -le 7; $col++) {
    $ws.Cells.Item(1, $col).Value2 = $newHeaders[$col - 1]
}
